$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-9 from 45243 (2023-11-13)
# to 45244 (2023-11-14), keeping existing date formatting/style.
foreach ($r in 2..9) {
    $ws.Cells.Item($r, 3).Value = 45244
}
